$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "ElasticSearch"
$ws.Range("A20").Value = "ReactNative"

$ws.Range("A21").Select()
